# This script inserts a new data record at row 24 of the "Mango" price
# sheet, pushing all existing records from row 24 onward down by one row
# (the last record, originally on row 135, ends up on new row 136).
#
# Only the columns that actually vary per record are shifted:
#   D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
#   O (Precio maximo), P (Precio promedio ponderado), R (Origen),
#   S (Precio $/Kg)
# All the other columns (A,B,C,E,F,G,H,I,J,K,Q,T) are identical on every
# data row, so there is nothing to shift there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 24
$lastDataRowBefore = 135
$lastDataRowAfter = 136

$shiftCols = @(4, 12, 13, 14, 15, 16, 18, 19)   # D, L, M, N, O, P, R, S
$allCols = 1..20                                 # A .. T

# First, create the brand new last row by duplicating every column from
# what was previously the last row (row 135), since row 136 did not
# exist before.
foreach ($c in $allCols) {
    $ws.Cells.Item($lastDataRowAfter, $c).Value = $ws.Cells.Item($lastDataRowBefore, $c).Value2
}
# The date column (D) carries a date number format; copy that format
# explicitly onto the newly created row so it keeps displaying as a date.
$ws.Cells.Item($lastDataRowAfter, 4).NumberFormat = $ws.Cells.Item($lastDataRowBefore, 4).NumberFormat

# Now shift the variable columns of the remaining existing rows down by
# one, starting from the bottom so source data isn't overwritten before
# it's copied. (Row 136's variable columns will be re-written here using
# row 135's data, which is the correct final value anyway.)
for ($r = $lastDataRowAfter; $r -gt $firstDataRow; $r--) {
    $srcRow = $r - 1
    foreach ($c in $shiftCols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($srcRow, $c).Value2
    }
}

# Write the brand-new record into row 24.
$ws.Cells.Item($firstDataRow, 4).Value  = 44749    # Fecha
$ws.Cells.Item($firstDataRow, 13).Value = 200       # Volumen
$ws.Cells.Item($firstDataRow, 14).Value = 8000      # Precio minimo
$ws.Cells.Item($firstDataRow, 15).Value = 8000      # Precio maximo
$ws.Cells.Item($firstDataRow, 16).Value = 8000      # Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 18).Value = "Brasil"  # Origen
$ws.Cells.Item($firstDataRow, 19).Value = 2000      # Precio $/Kg
# Calidad (L24) stays "Primera", unchanged.
